# Apply updated evidence / candidate-response columns (E, F, G, I, J, K, L)
# for the 5 data rows (rows 2-6), per the authoring commit "new changes
# created with all models" -- adds candidate_response1/candidate_response2
# data and replaces final_answer / original_response / final_ans_dict /
# candidate_response0 text with the new model outputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ques_id 0 (first animal to land on the moon)
$ws.Range("E2").Value = 'No animal has ever been sent to the moon.'
$ws.Range("F2").Value = 'test of new sc'
$ws.Range("G2").Value = 'No animal has ever been sent to the moon.'
$ws.Range("I2").Value = '{''QueryID'': ''0'', ''answer_box'': [{''name'': ''What was the first animal to land on the moon?'', ''url'': ''https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens.'', ''snippet'': ''Answer and Explanation: No animals were ever sent to the Moon. Although, since humans are technically animals, one could say that the first animal sent to the Moon was Neil Armstrong. He belonged to the species Homo sapiens.''}], ''related_questions'': [{''name'': ''What was the first animal to come on land?'', ''url'': ''https://bigthink.com/hard-science/first-creature-to-live-on-land/#:~:text=Scientists%20think%20an%20insect%20similar,first%2Dever%20land%2Ddweller.&text=An%20ancient%20millipede%2Dlike%20creature,creature%20to%20live%20on%20land.'', ''snippet'': ''Scientists think an insect similar to the modern millipede crawled around Scotland 425 million years ago, making it the first-ever land-dweller. An ancient millipede-like creature living in Scotland may have been the first creature to live on land.''}, {''name'': ''What was the name of the first land on the moon?'', ''url'': "https://www.rmg.co.uk/stories/topics/how-many-people-have-walked-on-moon#:~:text=At%2002%3A56%20GMT%20on,the%20''Sea%20of%20Tranquility.''", ''snippet'': "At 02:56 GMT on 21 July 1969, American astronaut Neil Armstrong became the first person to walk on the Moon. He stepped out of the Apollo 11 lunar module and onto the Moon''s surface, in an area called the ''Sea of Tranquility.''"}, {''name'': ''What were the first animals in space?'', ''url'': ''https://www.space.com/animals-in-space#:~:text=The%20first%20animals%20to%20reach,flights%20for%20about%20a%20decade.'', ''snippet'': ''The first animals to reach space were fruit flies that the United States launched aboard captured German rockets in 1947. The first mammal to reach space was a rhesus monkey named Albert II, who flew two years later. Both these missions were suborbital, as were all animal flights for about a decade.''}, {''name'': ''Was a dog the first on the moon?'', ''url'': ''https://www.quora.com/Why-did-they-send-a-dog-to-the-Moon-first-Why-didnt-they-choose-another-animal#:~:text=No%20dog%20has%20ever%20been,on%20the%20Soviet%20Sputnik%202.'', ''snippet'': ''No dog has ever been sent to the moon. The first (deliberate*) living thing launched into orbit was a dog name Laika, on the Soviet Sputnik 2.''}], ''organic_results'': [{''name'': ''What was the first animal in space?'', ''url'': ''https://www.rmg.co.uk/stories/topics/what-was-first-animal-space'', ''snippet'': ''The first animals to orbit the moon and return to Earth were two Russian tortoises aboard Zond 5. On 15 September 1968, the tortoises were launched with plants, ...''}, {''name'': ''Is it true that the first animal that went to the moon was a dog?'', ''url'': ''https://www.quora.com/Is-it-true-that-the-first-animal-that-went-to-the-moon-was-a-dog'', ''snippet'': ''No dog has ever been sent to the moon. The first (deliberate*) living thing launched into orbit was a dog name Laika, on the Soviet Sputnik 2.''}]}'
$ws.Range("J2").Value = '{''Answer:'': ''No animal has ever been sent to the moon.''}'
$ws.Range("K2").Value = '{''Answer:'': ''No animal has ever been sent to the moon.'', ''Source:'': ''https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens.'', ''Premise of the Question:'': ''The question assumes that animals have been sent to the moon, and asks which animal was the first to land on the moon. However, the answer provided in the source states that no animals have ever been sent to the moon, and therefore, there is no first animal to have landed on the moon.''}'
$ws.Range("L2").Value = '{''Answer:'': ''No animal has ever been sent to the moon.'', ''Source:'': ''https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens.'', ''Premise of the Question:'': ''The question assumes that animals have been sent to the moon, and asks which animal was the first to land on the moon. However, the answer provided in the source states that no animals have ever been sent to the moon, and therefore, there is no first animal to land on the moon.''}'

# Row 3 - ques_id 1 (Leonardo DiCaprio's third child)
$ws.Range("E3").Value = 'None.'
$ws.Range("F3").Value = 'test of new sc'
$ws.Range("G3").Value = 'None.'
$ws.Range("I3").Value = '{''QueryID'': ''1'', ''related_questions'': [{''name'': ''How many kids does Leo DiCaprio have?'', ''url'': ''https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have'', ''snippet'': ''He has never been married and has no children. How many relationships did Leonardo DiCaprio had in the past? Bridget Hall (1994): The model was the first woman linked to DiCaprio after the start of his movie career. Kristen Zang (1996): DiCaprio and Zang spent time together before his meteoric rise with Titanic.''}, {''name'': "What is Leonardo DiCaprio''s birth name?", ''url'': ''https://en.wikipedia.org/wiki/Leonardo_DiCaprio'', ''snippet'': ''Leonardo Wilhelm DiCaprio (/diˈkæprioʊ, dɪ-/; Italian: [diˈkaːprjo]; born November 11, 1974) is an American actor and film producer.''}, {''name'': "What is Leonardo DiCaprio''s longest relationship?", ''url'': ''https://www.yahoo.com/entertainment/leonardo-dicaprio-exes-said-dating-210507180.html'', ''snippet'': "One of Leonardo DiCaprio''s longest relationships was with Gisele Bundchen, whom he dated from 2001 to 2005. Bundchen has spoken about their relationship the most out of all of his exes. In a 2009 interview with Vanity Fair, she spoke about it, and had quite nice things to say."}, {''name'': ''How many siblings does Leonardo DiCaprio have?'', ''url'': ''https://kahawatungu.com/leonardo-dicaprio-siblings/'', ''snippet'': "DiCaprio''s parents divorced when he was a year old, and he was raised by his mother. He does not have any biological siblings, but he has a stepbrother named Adam Farrar. The relationship between DiCaprio and Farrar has been strained, and they have reportedly not been in contact for years."}], ''organic_results'': [{''name'': ''Leonardo DiCaprio'', ''url'': ''https://en.wikipedia.org/wiki/Leonardo_DiCaprio'', ''snippet'': "Leonardo Wilhelm DiCaprio was born on November 11, 1974, in Los Angeles, California. · DiCaprio''s parents named him Leonardo because his pregnant mother first ..."}, {''name'': "What Kate Winslet''s Children Call Leonardo DiCaprio", ''url'': ''https://www.nickiswift.com/455252/what-kate-winslets-children-call-leonardo-dicaprio/'', ''snippet'': "Leonardo DiCaprio and Kate Winslet''s unsinkable friendship ... These days, Kate Winslet is a proud mama of three children. She welcomed Mia Honey ..."}]}'
$ws.Range("J3").Value = '{''Answer:'': ''None.'', ''Source:'': ''https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have.'', ''Premise of the Question:'': ''The premise of this question is that the asker is querying about the number of children Leonardo DiCaprio has. However, according to the reference source, Leonardo DiCaprio does not have any children.''}'
$ws.Range("K3").Value = '{''Answer:'': ''None.'', ''Source:'': ''https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have.'', ''Premise of the Question:'': ''The premise of the question is false. Leonardo DiCaprio does not have any children.''}'
$ws.Range("L3").Value = '{''Answer:'': ''None.'', ''Source:'': ''https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have'', ''Premise of the Question:'': ''The premise of the question is false. Leonardo DiCaprio does not have any children, as stated in the source. Therefore, he does not have a third child.''}'

# Row 4 - ques_id 2 (first human to land on Mars)
$ws.Range("E4").Value = 'It seems a proper response could not be generated.'
$ws.Range("F4").Value = 'test of new sc'
$ws.Range("G4").Value = 'It seems a proper response could not be generated.'
$ws.Range("I4").Value = '{''QueryID'': ''2'', ''related_questions'': [{''name'': ''When did man first land on Mars?'', ''url'': "https://en.wikipedia.org/wiki/Mars_landing#:~:text=Soviet%20Union''s%20Mars%203%2C%20which,have%20conducted%20Mars%20landings%20successfully.", ''snippet'': "Soviet Union''s Mars 3, which landed in 1971, was the first successful Mars landing. As of 2023, the Soviet Union, United States and China have conducted Mars landings successfully."}, {''name'': ''What year will the first person be on Mars?'', ''url'': ''https://en.wikipedia.org/wiki/Human_mission_to_Mars#:~:text=The%20first%20crewed%20Mars%20Mission,is%20proposed%20for%20the%202030s.'', ''snippet'': ''The first crewed Mars Mission, which would include sending astronauts to Mars, orbiting Mars, and a return to Earth, is proposed for the 2030s.''}, {''name'': ''Who went 1st on Mars?'', ''url'': "https://en.wikipedia.org/wiki/Exploration_of_Mars#:~:text=Starting%20in%201960%2C%20the%20Soviets,1965%2C%20by%20NASA''s%20Mariner%204.", ''snippet'': "Starting in 1960, the Soviets launched a series of probes to Mars including the first intended flybys and hard (impact) landing (Mars 1962B). The first successful flyby of Mars was on 14–15 July 1965, by NASA''s Mariner 4."}, {''name'': ''When was the first US Mars landing?'', ''url'': ''https://www.planetary.org/space-missions/viking#:~:text=Viking%201%20launched%20on%20August,touched%20down%20on%20September%203.'', ''snippet'': ''Viking 1 launched on August 20, 1975. It arrived in Mars orbit on June 19, 1976 and the lander touched down on July 20, 1976. Viking 2 launched less than a month after Viking 1 on September 9, 1975. The spacecraft arrived in orbit on August 7, 1976 and the lander touched down on September 3.''}], ''organic_results'': [{''name'': ''Human mission to Mars'', ''url'': ''https://en.wikipedia.org/wiki/Human_mission_to_Mars'', ''snippet'': ''The first crewed Mars Mission, which would include sending astronauts to Mars, orbiting Mars, and a return to Earth, is proposed for the 2030s.''}, {''name'': ''Mars Exploration Timeline - the NSSDCA'', ''url'': ''https://nssdc.gsfc.nasa.gov/planetary/chronology_mars.html'', ''snippet'': ''Mission Timeline ; Mariner 6 - 25 February 1969 - Mars Flyby ; Mariner 7 - 27 March 1969 - Mars Flyby''}]}'
$ws.Range("J4").Value = '{''Answer:'': ''Not specified.'', ''Source:'': "https://en.wikipedia.org/wiki/Mars_landing#:~:text=Soviet%20Union''s%20Mars%203%2C%20which,have%20conducted%20Mars%20landings%20successfully.", ''Premise of the Question:'': ''The premise of the question is that a human has landed on Mars.The source provided does not specify a specific year for the first human landing on Mars. It only states that the Soviet Union\''s Mars 3 mission, which landed in 1971, was the first successful Mars landing, and that as of 2023, the Soviet Union, United States, and China have conducted Mars landings successfully.Therefore, the answer to the question "What year did the first human land on Mars?" cannot be determined based on the provided source.Let\''s think it step by step.1. The question asks about the first human landing on Mars.2. The provided source mentions the first successful Mars landing, which was the Soviet Union\''s Mars 3 mission in 1971.3. However, the source does not specify a specific year for the first human landing on Mars.4. Therefore, the answer to the question cannot be determined based on the provided source.Is there anything else I can help you with?''}'
$ws.Range("K4").Value = '{''Answer:'': ''No human has landed on Mars yet.''}'
$ws.Range("L4").Value = '{''Answer:'': ''No human has landed on Mars yet.''}'

# Row 5 - ques_id 3 (Zodiac killer)
$ws.Range("E5").Value = 'It seems a proper response could not be generated.'
$ws.Range("F5").Value = 'test of new sc'
$ws.Range("G5").Value = 'It seems a proper response could not be generated.'
$ws.Range("I5").Value = '{''QueryID'': ''3'', ''answer_box'': [{''name'': ''Zodiac Killer: Biography, Serial Killer, Criminal'', ''url'': ''https://www.biography.com/crime/zodiac-killer#:~:text=True%2Dcrime%20author%20and%20former,as%20the%20most%20likely%20suspect.'', ''snippet'': "True-crime author and former San Francisco Chronicle cartoonist Robert Graysmith wrote two separate works on the killer (1986''s Zodiac and 2002''s Zodiac Unmasked), ultimately identifying a man named Arthur Leigh Allen as the most likely suspect."}], ''related_questions'': [{''name'': "What is the Zodiac Killer''s real name?", ''url'': ''https://en.wikipedia.org/wiki/Zodiac_Killer#:~:text=Of%20the%20four%20ciphers%20he,offender%20who%20died%20in%201992.'', ''snippet'': ''Of the four ciphers he produced, two remain unsolved, and one was cracked only in 2020. While many theories regarding the identity of the killer have been suggested, the only suspect authorities ever publicly named was Arthur Leigh Allen, a former elementary school teacher and convicted sex offender who died in 1992.''}, {''name'': ''Is Arthur Leigh Allen the Zodiac?'', ''url'': ''https://screenrant.com/all-evidence-arthur-leigh-allen-not-zodiac-killer/#:~:text=Arthur%20Leigh%20Allen%20was%20the,the%20case%20is%20still%20open.'', ''snippet'': ''Arthur Leigh Allen was the most likely suspect in the Zodiac case, but there was not enough physical evidence to prove his guilt. The circumstantial evidence and fear in the community led to Allen being falsely identified as the killer. The identity of the Zodiac Killer remains a mystery, and the case is still open.''}, {''name'': ''Were there 2 Zodiac killers?'', ''url'': ''https://www.distractify.com/p/how-many-zodiac-killers-were-there#:~:text=When%20Choon%2DJae%20was%20later,anything%20is%20possible.'', ''snippet'': ''When Choon-Jae was later caught in 2018, he could no longer be compared to the one and only Zodiac. So there are at least three other killers who have been compared to the Zodiac Killer, and while experts agree that the original Zodiac Killer is likely just one man … anything is possible.''}, {''name'': ''Did they ever catch Jack the Ripper?'', ''url'': ''https://www.britannica.com/biography/Jack-the-Ripper#:~:text=Jack%20the%20Ripper%20was%20an,macabre%20tourist%20industry%20in%20London.'', ''snippet'': "Jack the Ripper was an English serial killer. Between August and November 1888, he murdered at least five women—all prostitutes—in or near the Whitechapel district of London''s East End. Jack the Ripper was never identified or arrested. Today the murder sites are the locus of a macabre tourist industry in London."}], ''organic_results'': [{''name'': ''Zodiac Killer'', ''url'': ''https://en.wikipedia.org/wiki/Zodiac_Killer'', ''snippet'': ''For the Japanese Zodiac copycat, see Kobe child murders. The Zodiac Killer is the pseudonym of an unidentified serial killer who operated in Northern California ...''}, {''name'': "Zodiac Killer''s identity has been reportedly revealed, DNA ...", ''url'': ''https://www.marca.com/en/lifestyle/us-news/2023/05/20/64680e6de2704e873f8b4591.html'', ''snippet'': "The Zodiac Killer''s identity has been reportedly found by the FBI, they believe that although the infamous murderer that was active during ..."}]}'
$ws.Range("J5").Value = '{''Answer:'': ''The name of the Zodiac killer is Arthur Leigh Allen.''}'
$ws.Range("K5").Value = '{''Answer:'': ''The name of the Zodiac killer is Arthur Leigh Allen.''}'
$ws.Range("L5").Value = '{''Answer:'': ''Arthur Leigh Allen'', ''Source:'': ''https://en.wikipedia.org/wiki/Zodiac_Killer'', ''Premise of the Question:'': "The Zodiac Killer''s real name was Arthur Leigh Allen, a former elementary school teacher and convicted sex offender who died in 1992. This is the most likely suspect identified by authorities, although the case remains open and the identity of the killer is still a mystery."}'

# Row 6 - ques_id 4 (P vs NP)
$ws.Range("E6").Value = 'It seems a proper response could not be generated.'
$ws.Range("F6").Value = 'test of new sc'
$ws.Range("G6").Value = 'It seems a proper response could not be generated.'
$ws.Range("I6").Value = '{''QueryID'': ''4'', ''related_questions'': [{''name'': ''Are all problems solvable in polynomial time?'', ''url'': ''https://www.cs.umd.edu/class/fall2017/cmsc451-0101/Handouts/sol5.pdf'', ''snippet'': ''(i) All NP-complete problems are solvable in polynomial time: Yes. Every problem in NP is polynomially reducible to SAT, and SAT is reducible to every NP-hard problem. Therefore, a polynomial time solution to any NP-hard problem (such as 3Col) implies that every problem in NP can be solved in polynomial time.''}, {''name'': ''Why can we solve every problem in NP in exponential time?'', ''url'': ''https://www.quora.com/Can-all-NP-problems-be-solved-in-exponential-time#:~:text=Yes%2C%20cause%20the%20class%20of,in%20that%20computational%20complexity%20class).'', ''snippet'': ''Yes, cause the class of NP is contained within the class EXP (EXPTIME), while EXP is the class of all problems solvable (decidable) within exponential running-time (relative to the representation of the size of the non-unary input to the problem in that computational complexity class).''}, {''name'': ''What does it mean to be solvable in polynomial time?'', ''url'': ''https://www.britannica.com/science/polynomial-time-algorithm#:~:text=computational%20problems&text=%E2%80%A6can%20be%20solved%20in%20%E2%80%9Cpolynomial,the%20input%20for%20the%20problem.'', ''snippet'': ''computational problems …can be solved in “polynomial time,” which means that an algorithm exists for its solution such that the number of steps in the algorithm is bounded by a polynomial function of n, where n corresponds to the length of the input for the problem.''}, {''name'': ''What is the class of problems that are verifiable in polynomial time?'', ''url'': ''https://en.wikipedia.org/wiki/P_versus_NP_problem#:~:text=The%20class%20of%20questions%20where,for%20%22nondeterministic%20polynomial%20time%22.'', ''snippet'': ''The class of questions where an answer can be verified in polynomial time is NP, standing for "nondeterministic polynomial time".''}], ''organic_results'': [{''name'': ''P versus NP problem'', ''url'': ''https://en.wikipedia.org/wiki/P_versus_NP_problem'', ''snippet'': ''Informally, it asks whether every problem whose solution can be quickly verified can also be quickly solved. Here, quickly means an algorithm that solves ...''}, {''name'': ''Is there a task that is solvable in polynomial time but not ...'', ''url'': ''https://www.quora.com/Is-there-a-task-that-is-solvable-in-polynomial-time-but-not-verifiable-in-polynomial-time'', ''snippet'': ''For some problems, the answer can be verified to be correct in Polynomial Time, even if there is no known way of solving the original problem in ...''}]}'
$ws.Range("J6").Value = '{''Answer:'': ''No, not all quickly verifiable problems are quickly solvable.''}'
$ws.Range("K6").Value = '{''Answer:'': ''No, not all quickly verifiable problems are quickly solvable.''}'
$ws.Range("L6").Value = '{''Answer:'': ''No, not all quickly verifiable problems are quickly solvable.''}'
